$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 962.5
$ws.Range("J17").Value = 1004.5455
$ws.Range("L17").Value = 3013.6365
$ws.Range("N17").Value = -3349.6365
$ws.Range("H48").Value = 4749.8335
$ws.Range("I48").Value = 3000
$ws.Range("J48").Value = 5099.8
$ws.Range("K48").Value = 9000
$ws.Range("L48").Value = 15299.4
$ws.Range("M48").Value = -8708
$ws.Range("N48").Value = -15883.4
$ws.Range("H56").Value = 4749.8335
$ws.Range("I56").Value = 3000
$ws.Range("J56").Value = 5099.8
$ws.Range("K56").Value = 9000
$ws.Range("L56").Value = 15299.4
$ws.Range("M56").Value = -8466
$ws.Range("N56").Value = -16367.4
$ws.Range("H86").Value = 3637
$ws.Range("I86").Value = 4967.1665
$ws.Range("J86").Value = 2040.8
$ws.Range("K86").Value = 4967.1665
$ws.Range("L86").Value = 2040.8
$ws.Range("M86").Value = -3844.1665
$ws.Range("N86").Value = -4286.8
$ws.Range("H89").Value = 3637
$ws.Range("I89").Value = 4967.1665
$ws.Range("J89").Value = 2040.8
$ws.Range("K89").Value = 24835.8325
$ws.Range("L89").Value = 10204
$ws.Range("M89").Value = -19219.8325
$ws.Range("N89").Value = -21436
$ws.Range("H96").Value = 1872.5454
$ws.Range("J96").Value = 977.25
$ws.Range("L96").Value = 2931.75
$ws.Range("N96").Value = -5677.75
$ws.Range("H116").Value = 2401.4849
$ws.Range("I116").Value = 2686.5833
$ws.Range("J116").Value = 2238.5715
$ws.Range("K116").Value = 2686.5833
$ws.Range("L116").Value = 2238.5715
$ws.Range("M116").Value = 755.4167000000002
$ws.Range("N116").Value = -9122.5715
$ws.Range("H138").Value = 2733.2812
$ws.Range("I138").Value = 2743.6667
$ws.Range("J138").Value = 2732.2068
$ws.Range("K138").Value = 8231.000100000001
$ws.Range("L138").Value = 8196.6204
$ws.Range("M138").Value = -3091.000100000001
$ws.Range("N138").Value = -18476.6204

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10822.521
$ws.Range("I32").Value = 8076.5845
$ws.Range("J32").Value = 17433.111
$ws.Range("K32").Value = 8076.5845
$ws.Range("L32").Value = 17433.111
$ws.Range("M32").Value = -7789.5845
$ws.Range("N32").Value = -18007.111
$ws.Range("H45").Value = 1123.875
$ws.Range("I45").Value = 975.75
$ws.Range("J45").Value = 1272
$ws.Range("K45").Value = 975.75
$ws.Range("L45").Value = 1272
$ws.Range("M45").Value = -598.75
$ws.Range("N45").Value = -2026
$ws.Range("H63").Value = 2027
$ws.Range("I63").Value = 1874.4546
$ws.Range("J63").Value = 2506.4285
$ws.Range("K63").Value = 1874.4546
$ws.Range("L63").Value = 2506.4285
$ws.Range("M63").Value = -1188.4546
$ws.Range("N63").Value = -3878.4285
$ws.Range("H66").Value = 2027
$ws.Range("I66").Value = 1874.4546
$ws.Range("J66").Value = 2506.4285
$ws.Range("K66").Value = 9372.273000000001
$ws.Range("L66").Value = 12532.1425
$ws.Range("M66").Value = -5940.273000000001
$ws.Range("N66").Value = -19396.1425
$ws.Range("H74").Value = 2104.1155
$ws.Range("I74").Value = 1149.7778
$ws.Range("J74").Value = 4251.375
$ws.Range("K74").Value = 1149.7778
$ws.Range("L74").Value = 4251.375
$ws.Range("M74").Value = -275.7778000000001
$ws.Range("N74").Value = -5999.375
$ws.Range("H77").Value = 2104.1155
$ws.Range("I77").Value = 1149.7778
$ws.Range("J77").Value = 4251.375
$ws.Range("K77").Value = 5748.889
$ws.Range("L77").Value = 21256.875
$ws.Range("M77").Value = -1380.889
$ws.Range("N77").Value = -29992.875
$ws.Range("H102").Value = 27782330
$ws.Range("I102").Value = 41668496
$ws.Range("K102").Value = 41668496
$ws.Range("M102").Value = -41666874
$ws.Range("H132").Value = 2645.745
$ws.Range("I132").Value = 2022.1389
$ws.Range("J132").Value = 4142.4
$ws.Range("K132").Value = 6066.4167
$ws.Range("L132").Value = 12427.2
$ws.Range("M132").Value = -3536.4167
$ws.Range("N132").Value = -17487.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 340.5
$ws.Range("I12").Value = 340.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 340.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -172.5
$ws.Range("N12").ClearContents()
$ws.Range("H134").Value = 5352.7085
$ws.Range("I134").Value = 1005.1429
$ws.Range("J134").Value = 11439.3
$ws.Range("K134").Value = 3015.4287
$ws.Range("L134").Value = 34317.89999999999
$ws.Range("M134").Value = -480.4287000000004
$ws.Range("N134").Value = -39387.89999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1370.04
$ws.Range("I31").Value = 1370.04
$ws.Range("K31").Value = 1370.04
$ws.Range("M31").Value = -1075.04
$ws.Range("H34").Value = 1370.04
$ws.Range("I34").Value = 1370.04
$ws.Range("K34").Value = 1370.04
$ws.Range("M34").Value = -1168.04
$ws.Range("H97").Value = 29800
$ws.Range("J97").Value = 29800
$ws.Range("L97").Value = 29800
$ws.Range("N97").Value = -31782
$ws.Range("H99").Value = 1543.6666
$ws.Range("I99").Value = 1494.1428
$ws.Range("K99").Value = 1494.1428
$ws.Range("M99").Value = 3.857199999999921
$ws.Range("H126").Value = 1543.6666
$ws.Range("I126").Value = 1494.1428
$ws.Range("K126").Value = 4482.428400000001
$ws.Range("M126").Value = -2012.428400000001
$ws.Range("H131").Value = 19999
$ws.Range("J131").Value = 19999
$ws.Range("L131").Value = 19999
$ws.Range("N131").Value = -30079
$ws.Range("H132").Value = 1561.9333
$ws.Range("I132").Value = 1212.125
$ws.Range("K132").Value = 3636.375
$ws.Range("M132").Value = -1106.375
$ws.Range("H134").Value = 11112507
$ws.Range("I134").Value = 1384.5
$ws.Range("J134").Value = 38462964
$ws.Range("K134").Value = 4153.5
$ws.Range("L134").Value = 115388892
$ws.Range("M134").Value = -1618.5
$ws.Range("N134").Value = -115393962

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 17362218
$ws.Range("I129").Value = 37037556
$ws.Range("J129").Value = 5557017.5
$ws.Range("K129").Value = 111112668
$ws.Range("L129").Value = 16671052.5
$ws.Range("M129").Value = -111107668
$ws.Range("N129").Value = -16681052.5
$ws.Range("H131").Value = 22762246
$ws.Range("I131").Value = 71428970
$ws.Range("J131").Value = 51110.3
$ws.Range("K131").Value = 214286910
$ws.Range("L131").Value = 153330.9
$ws.Range("M131").Value = -214281870
$ws.Range("N131").Value = -163410.9
$ws.Range("H137").Value = 34098772
$ws.Range("I137").Value = 83336210
$ws.Range("J137").Value = 11320.385
$ws.Range("K137").Value = 250008630
$ws.Range("L137").Value = 33961.155
$ws.Range("M137").Value = -250003530
$ws.Range("N137").Value = -44161.155

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4050.3635
$ws.Range("I80").Value = 2004
$ws.Range("K80").Value = 2004
$ws.Range("M80").Value = -1006
$ws.Range("H83").Value = 4050.3635
$ws.Range("I83").Value = 2004
$ws.Range("K83").Value = 10020
$ws.Range("M83").Value = -5028
$ws.Range("H102").Value = 2726.4644
$ws.Range("I102").Value = 1944.3889
$ws.Range("K102").Value = 1944.3889
$ws.Range("M102").Value = -322.3888999999999
$ws.Range("H123").Value = 24260.8
$ws.Range("J123").Value = 24260.8
$ws.Range("L123").Value = 24260.8
$ws.Range("N123").Value = -29160.8
$ws.Range("H126").Value = 2188.75
$ws.Range("I126").Value = 1922
$ws.Range("J126").Value = 2633.3333
$ws.Range("K126").Value = 5766
$ws.Range("L126").Value = 7899.999899999999
$ws.Range("M126").Value = -3296
$ws.Range("N126").Value = -12839.9999
$ws.Range("H132").Value = 6924.269
$ws.Range("I132").Value = 10276.308
$ws.Range("J132").Value = 3572.2307
$ws.Range("K132").Value = 30828.924
$ws.Range("L132").Value = 10716.6921
$ws.Range("M132").Value = -28298.924
$ws.Range("N132").Value = -15776.6921

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2062.5715
$ws.Range("I7").Value = 1587.6
$ws.Range("J7").Value = 3250
$ws.Range("K7").Value = 1587.6
$ws.Range("L7").Value = 3250
$ws.Range("M7").Value = -1475.6
$ws.Range("N7").Value = -3474
$ws.Range("H40").Value = 3594.75
$ws.Range("I40").Value = 2356.5454
$ws.Range("K40").Value = 2356.5454
$ws.Range("M40").Value = -2220.5454
$ws.Range("H68").Value = 1318.4166
$ws.Range("I68").Value = 1318.4166
$ws.Range("K68").Value = 1318.4166
$ws.Range("M68").Value = -569.4166
$ws.Range("H71").Value = 1318.4166
$ws.Range("I71").Value = 1318.4166
$ws.Range("K71").Value = 6592.083000000001
$ws.Range("M71").Value = -2848.083000000001
$ws.Range("H122").Value = 15627656
$ws.Range("I122").Value = 17859678
$ws.Range("K122").Value = 53579034
$ws.Range("M122").Value = -53576584
$ws.Range("H126").Value = 2062.5715
$ws.Range("I126").Value = 1587.6
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 4762.799999999999
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -2292.799999999999
$ws.Range("N126").Value = -14690
$ws.Range("H132").Value = 2440.1562
$ws.Range("I132").Value = 1999.35
$ws.Range("J132").Value = 3174.8333
$ws.Range("K132").Value = 5998.049999999999
$ws.Range("L132").Value = 9524.499899999999
$ws.Range("M132").Value = -3468.049999999999
$ws.Range("N132").Value = -14584.4999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1823.6538
$ws.Range("J81").Value = 1916.2858
$ws.Range("L81").Value = 3832.5716
$ws.Range("N81").Value = -5954.5716
$ws.Range("H84").Value = 1823.6538
$ws.Range("J84").Value = 1916.2858
$ws.Range("L84").Value = 19162.858
$ws.Range("N84").Value = -29770.858
$ws.Range("H96").Value = 957.1539
$ws.Range("I96").Value = 1460.6
$ws.Range("J96").Value = 642.5
$ws.Range("K96").Value = 1460.6
$ws.Range("L96").Value = 642.5
$ws.Range("M96").Value = -87.59999999999991
$ws.Range("N96").Value = -3388.5
$ws.Range("H122").Value = 13891604
$ws.Range("I122").Value = 14708581
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 44125743
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -44123293
$ws.Range("N122").Value = -13900
$ws.Range("H136").Value = 1678.55
$ws.Range("J136").Value = 2578.5715
$ws.Range("L136").Value = 7735.7145
$ws.Range("N136").Value = -12835.7145

Write-Host "Applied 266 cell updates across 8 sheets."
